# "all single filter scripts in CTDC"
# Rewrites the CasesTab query (B2), and adds a new FilesTab row (row 3)
# with the file-level query (B3) and an updated stat query (C3), mirroring
# the C2 stat-query text (with a 1-space-wider indent variant).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New text for B2 (CasesTab query), rewritten ----
$casesTabQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
  WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# ---- Updated text for C2 (StatQuery, first variant) ----
$statQueryV1 = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
       WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# ---- Text for C3 (StatQuery, second variant - extra indent before WHERE) ----
$statQueryV2 = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE c.ethnicity = "HISPANIC_OR_LATINO"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# ---- Text for B3 (FilesTab query) ----
$filesTabQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE c.ethnicity = "HISPANIC_OR_LATINO"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Here-strings keep a trailing newline before the closing '@ - strip it so
# the cell text matches the source exactly (no trailing blank line).
$casesTabQuery = $casesTabQuery.TrimEnd("`r","`n")
$statQueryV1   = $statQueryV1.TrimEnd("`r","`n")
$statQueryV2   = $statQueryV2.TrimEnd("`r","`n")
$filesTabQuery = $filesTabQuery.TrimEnd("`r","`n")

# ---- Row 2: update CasesTab's query (B2) and the StatQuery text (C2) ----
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQueryV1

# ---- Row 3: brand-new FilesTab row ----
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("C3").Value = $statQueryV2
$ws.Range("D3").Value = "TC01_Trials_Filter_Ethnicity-HispLatino_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_Ethnicity-HispLatino_WebData.xlsx"

$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# ---- Row heights (auto-sized by Excel for the new wrapped content) ----
$ws.Rows(2).RowHeight = 188.5
$ws.Rows(3).RowHeight = 409.5

# ---- View state: scrolled to row 3, zoomed to 70%, selection on C3 ----
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("C3").Select()
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 1
